$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of row 39 (a previously-added "finished" entry) onto
# the new row 41 so fonts/fills/number formats/borders all match exactly.
$ws.Range("A39:H39").Copy() | Out-Null
$ws.Range("A41:H41").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(41).RowHeight = 56

$ws.Range("A41").Value = "309. Best Time to Buy and Sell Stock with Cooldown"
$ws.Range("B41").Value = "Medium"
$ws.Range("C41").Value = "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-with-cooldown/"
$ws.Range("D41").Value = (Get-Date -Year 2021 -Month 12 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E41").Value = "股票交易(状态机）"
$ws.Range("F41").Value = "涉及冷却时间的股票交易要用状态机"
$ws.Range("G41").Value = "未复习"
$ws.Range("H41").Value = "⭕"

$ws.Hyperlinks.Add(
    $ws.Cells.Item(41, 3),
    "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-with-cooldown/"
) | Out-Null

# Adding the hyperlink re-applies Excel's builtin "Hyperlink" style; restore
# the worksheet's own themed formatting (matching the other rows) afterwards.
$ws.Range("C39").Copy() | Out-Null
$ws.Range("C41").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A25").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("F48").Select() | Out-Null

$wb.Save()
